$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("my_excel")

$ws.Range("D2").Value = "10 394,78 kr"
$ws.Range("E2").Value = "330,32 kr"
$ws.Range("D3").Value = "110,22 kr"
$ws.Range("D4").Value = "100,00 kr"
$ws.Range("E5").Value = "0,00 kr"
$ws.Range("D6").Value = "0,00 kr"
$ws.Range("E6").Value = "0,00 kr"
